$wb = $excel.ActiveWorkbook

# Sheet 1 ("想去人数" / F column updates)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 17
$ws.Cells.Item(6, 6).Value = 1617
$ws.Cells.Item(7, 6).Value = 25
$ws.Cells.Item(9, 6).Value = 754
$ws.Cells.Item(10, 6).Value = 2729
$ws.Cells.Item(11, 6).Value = 2729
$ws.Cells.Item(12, 6).Value = 22
$ws.Cells.Item(13, 6).Value = 1843
$ws.Cells.Item(14, 6).Value = 626
$ws.Cells.Item(15, 6).Value = 313
$ws.Cells.Item(16, 6).Value = 719
$ws.Cells.Item(17, 6).Value = 12
$ws.Cells.Item(18, 6).Value = 6341
$ws.Cells.Item(19, 6).Value = 248
$ws.Cells.Item(20, 6).Value = 94
$ws.Cells.Item(22, 6).Value = 4416
$ws.Cells.Item(23, 6).Value = 892
$ws.Cells.Item(27, 6).Value = 2485
$ws.Cells.Item(36, 6).Value = 89
$ws.Cells.Item(37, 6).Value = 35
$ws.Cells.Item(38, 6).Value = 61
$ws.Cells.Item(39, 6).Value = 1532
$ws.Cells.Item(40, 6).Value = 37
$ws.Cells.Item(41, 6).Value = 1490

# Sheet 2 ("想去人数" / F column updates)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 185
$ws.Cells.Item(17, 6).Value = 150
$ws.Cells.Item(18, 6).Value = 343
$ws.Cells.Item(19, 6).Value = 274
$ws.Cells.Item(20, 6).Value = 526

# Sheet 3 ("想去人数" / F column updates)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 945
$ws.Cells.Item(8, 6).Value = 16

# Sheet 4 ("想去人数" / F column updates)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 945
$ws.Cells.Item(14, 6).Value = 17
$ws.Cells.Item(19, 6).Value = 2729
$ws.Cells.Item(20, 6).Value = 16
$ws.Cells.Item(22, 6).Value = 22
$ws.Cells.Item(23, 6).Value = 185
$ws.Cells.Item(24, 6).Value = 626
$ws.Cells.Item(25, 6).Value = 313
$ws.Cells.Item(26, 6).Value = 719
$ws.Cells.Item(27, 6).Value = 6341
$ws.Cells.Item(28, 6).Value = 248
$ws.Cells.Item(29, 6).Value = 94
$ws.Cells.Item(34, 6).Value = 2485
$ws.Cells.Item(39, 6).Value = 150
$ws.Cells.Item(40, 6).Value = 343
$ws.Cells.Item(41, 6).Value = 274
$ws.Cells.Item(42, 6).Value = 526
$ws.Cells.Item(44, 6).Value = 89
$ws.Cells.Item(45, 6).Value = 35
$ws.Cells.Item(46, 6).Value = 61
$ws.Cells.Item(48, 6).Value = 1532
$ws.Cells.Item(49, 6).Value = 37
